# Refresh cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
# Each row in the sheet holds one coin: columns D (Price) and E (Volume 1h change)
# are refreshed with the latest scraped values. Numeric-looking price strings are
# written with a leading apostrophe so Excel keeps them as text (matching the sheet's
# existing "Price" column, which stores values like "42.831.26" that are not valid
# numbers) and the style is reset to Normal right after so no extra number format
# sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '42.831.26'
$ws.Range("E2").Value = '  +0.58%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.281.88'
$ws.Range("E3").Value = '  -0.44%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.17%  '

# Row 5: BNB
$ws.Range("D5").Value = "'310.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.71%  '

# Row 6: Solana
$ws.Range("D6").Value = "'102.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.35%  '

# Row 7: XRP
$ws.Range("D7").Value = "'0.615"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.34%  '

# Row 8: USDC
$ws.Range("E8").Value = '  -0.11%  '

# Row 9: Cardano
$ws.Range("E9").Value = '  -0.97%  '

# Row 10: Avalanche
$ws.Range("D10").Value = "'38.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.94%  '

# Row 11: Dogecoin
$ws.Range("D11").Value = "'0.0899"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.98%  '

# Row 12: Polkadot
$ws.Range("D12").Value = "'8.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.22%  '

# Row 13: TRON
$ws.Range("D13").Value = "'0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.19%  '

# Row 14: Polygon
$ws.Range("E14").Value = '  +0.92%  '

# Row 15: Chainlink
$ws.Range("D15").Value = "'15.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.34%  '

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("D16").Value = '2.630.65'
$ws.Range("E16").Value = '  -0.37%  '

# Row 17: WrappedEther
$ws.Range("D17").Value = '2.276.57'
$ws.Range("E17").Value = '  -0.46%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '42.444.93'
$ws.Range("E18").Value = '  +0.13%  '

# Row 19: Uniswap
$ws.Range("D19").Value = "'7.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.46%  '

# Row 20: ShibaInu
$ws.Range("E20").Value = '  -1.39%  '

# Row 21: InternetComputer(DFINITY)
$ws.Range("D21").Value = "'13.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.59%  '

# Row 22: Litecoin
$ws.Range("D22").Value = "'73.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.41%  '

# Row 23: BitcoinCash
$ws.Range("D23").Value = "'266.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.20%  '

# Row 24: PancakeSwap
$ws.Range("D24").Value = "'3.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.40%  '

# Row 25: ImmutableX
$ws.Range("D25").Value = "'2.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.50%  '

# Row 27: Cosmos
$ws.Range("D27").Value = "'10.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.37%  '

# Row 28: Filecoin
$ws.Range("D28").Value = "'7.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +15.20%  '

# Row 29: Toncoin
$ws.Range("D29").Value = "'2.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.27%  '

# Row 30: EthereumClassic
$ws.Range("D30").Value = "'22.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.37%  '

# Row 31: InjectiveProtocol
$ws.Range("D31").Value = "'35.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.29%  '

# Row 32: Monero
$ws.Range("D32").Value = "'164.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.82%  '

# Row 33: Hedera
$ws.Range("D33").Value = "'0.0848"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.86%  '

# Row 34: Stellar
$ws.Range("E34").Value = '  -1.90%  '

# Row 35: WEMIXToken
$ws.Range("D35").Value = "'2.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.95%  '

# Row 36: Kaspa
$ws.Range("D36").Value = "'0.111"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.79%  '

# Row 37: RenderToken
$ws.Range("D37").Value = "'4.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.92%  '

# Row 38: VeChain
$ws.Range("E38").Value = '  -2.55%  '

# Row 39: LidoDAOToken
$ws.Range("E39").Value = '  +0.81%  '

# Row 40: NEARProtocol
$ws.Range("E40").Value = '  -4.53%  '

# Row 41: BitcoinSV
$ws.Range("D41").Value = "'107.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +12.28%  '

# Row 42: ARBITRUM
$ws.Range("D42").Value = "'1.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.19%  '

# Row 43: MultiversX
$ws.Range("D43").Value = "'71.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.64%  '

# Row 44: Algorand -> FirstDigitalUSD
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.06%  '

# Row 45: FirstDigitalUSD -> Algorand
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = "'0.225"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.03%  '

# Row 46: Celestia
$ws.Range("D46").Value = "'11.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.23%  '

# Row 47: Maker
$ws.Range("D47").Value = '1.720.62'
$ws.Range("E47").Value = '  +8.74%  '

# Row 48: Aave
$ws.Range("D48").Value = "'110.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.30%  '

# Row 49: ordi
$ws.Range("D49").Value = "'76.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.62%  '

# Row 50: THORChain
$ws.Range("D50").Value = "'5.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.54%  '

# Row 51: FraxShare
$ws.Range("D51").Value = "'8.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.30%  '
